# Rename the inline picture shapes that live in the document's header and
# footer stories:
#   - the Pearson logo (PNG) pictures:  image2.png -> image1.png
#   - the BTec logo (JPG) pictures:     image1.jpg -> image2.jpg
#
# wdPrimaryHeaderStory=7, wdPrimaryFooterStory=9,
# wdFirstPageHeaderStory=10, wdFirstPageFooterStory=11
$wdPrimaryHeaderStory    = 7
$wdPrimaryFooterStory    = 9
$wdFirstPageHeaderStory  = 10
$wdFirstPageFooterStory  = 11

$d = $word.ActiveDocument

foreach ($story in $d.StoryRanges) {
    if ($story.InlineShapes.Count -le 0) {
        continue
    }

    $storyType = $story.StoryType
    $isHeader = ($storyType -eq $wdPrimaryHeaderStory) -or ($storyType -eq $wdFirstPageHeaderStory)
    $isFooter = ($storyType -eq $wdPrimaryFooterStory) -or ($storyType -eq $wdFirstPageFooterStory)

    if (-not ($isHeader -or $isFooter)) {
        continue
    }

    $inlineShape = $story.InlineShapes.Item(1)

    # Selecting the shape first and renaming it through the resulting
    # Selection's InlineShapes collection is what actually commits the new
    # name back to the saved package for header/footer stories.
    $inlineShape.Select()
    $selectedShape = $word.Selection.InlineShapes.Item(1)

    if ($isHeader) {
        $selectedShape.Name = "image2.jpg"
    } else {
        $selectedShape.Name = "image1.png"
    }
}
